# Update "想去人数" (want-to-go count) figures to reflect newly generated output.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row 3 -> 266 to 269, row 4 -> 919 to 922
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 269
$wsExhibit.Range("F4").Value = 922

# Sheet "全部类型" (sheet4): row 4 -> 266 to 269, row 5 -> 919 to 922
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 269
$wsAll.Range("F5").Value = 922
